$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B7").Value = "asfd"
$ws.Range("C7").Value = "asdf"
$ws.Range("D7").Value = "asf"
$ws.Range("C8").Value = "ASDF"
$ws.Range("B9").Value = "asdf"
$ws.Range("D9").Value = "asdf"
$ws.Range("B10").Value = "asdf"
$ws.Range("C10").Value = "asdf"
$ws.Range("D10").Value = "asdf"
$ws.Range("J10").Select()
